$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Richness")
$ws.Range("A2").Value = 'healthy vs Grouppre_ltx'
$ws.Range("B2").Value = -53.6418918918919
$ws.Range("C2").Value = 11.7419507423117
$ws.Range("D2").Value = -4.5683969443506
$ws.Range("E2").Value = 0.000010737493714208
$ws.Range("F2").Value = 0.0000966374434278722
$ws.Range("A3").Value = 'healthy , pre_ltx - CZ vs NO'
$ws.Range("B3").Value = 11.2192192192192
$ws.Range("C3").Value = 10.4879627740154
$ws.Range("D3").Value = 1.06972340205245
$ws.Range("E3").Value = 0.286597899080415
$ws.Range("F3").Value = 0.360611710559031
$ws.Range("A4").Value = 'healthy vs Grouppre_ltx:CountryNO'
$ws.Range("B4").Value = 6.31338947643298
$ws.Range("D4").Value = 0.409876055554212
$ws.Range("E4").Value = 0.682527868900786
$ws.Range("F4").Value = 0.682527868900786
$ws.Range("F6").Value = 0.249171048977547
$ws.Range("A8").Value = 'healthy vs Grouppost_ltx'
$ws.Range("B8").Value = -28.0193428722841
$ws.Range("C8").Value = 9.23672572957162
$ws.Range("D8").Value = -3.03347135041365
$ws.Range("E8").Value = 0.00272705514175019
$ws.Range("F8").Value = 0.0122717481378758
$ws.Range("A9").Value = 'healthy , post_ltx - CZ vs NO'
$ws.Range("B9").Value = 11.2192192192191
$ws.Range("C9").Value = 11.2673273879844
$ws.Range("D9").Value = 0.995730294584617
$ws.Range("E9").Value = 0.320543742719139
$ws.Range("F9").Value = 0.360611710559031
$ws.Range("A10").Value = 'healthy vs Grouppost_ltx:CountryNO'
$ws.Range("B10").Value = -29.6195460166047
$ws.Range("D10").Value = -2.02470863408955
$ws.Range("E10").Value = 0.0441817364986902
$ws.Range("F10").Value = 0.0795271256976424

$ws = $wb.Worksheets.Item("Shannon")
$ws.Range("A2").Value = 'healthy vs Grouppre_ltx'
$ws.Range("B2").Value = -0.505199590305766
$ws.Range("C2").Value = 0.15931319845803
$ws.Range("D2").Value = -3.17110945731755
$ws.Range("E2").Value = 0.00186876305443549
$ws.Range("F2").Value = 0.0168188674899194
$ws.Range("G2").Value = '*'
$ws.Range("A3").Value = 'healthy , pre_ltx - CZ vs NO'
$ws.Range("B3").Value = 0.00655261583069825
$ws.Range("C3").Value = 0.142299259425117
$ws.Range("D3").Value = 0.046048137264878
$ws.Range("E3").Value = 0.963337973541263
$ws.Range("F3").Value = 0.963337973541263
$ws.Range("A4").Value = 'healthy vs Grouppre_ltx:CountryNO'
$ws.Range("B4").Value = 0.143633850099741
$ws.Range("D4").Value = 0.687282467598126
$ws.Range("E4").Value = 0.493049877549973
$ws.Range("F4").Value = 0.633921271135679
$ws.Range("F5").Value = 0.048753308414039
$ws.Range("F6").Value = 0.541565633906018
$ws.Range("F7").Value = 0.048753308414039
$ws.Range("A8").Value = 'healthy vs Grouppost_ltx'
$ws.Range("B8").Value = -0.133631433207174
$ws.Range("D8").Value = -1.20046961739458
$ws.Range("E8").Value = 0.231328990440202
$ws.Range("F8").Value = 0.416392182792363
$ws.Range("A9").Value = 'healthy , post_ltx - CZ vs NO'
$ws.Range("B9").Value = 0.006552615830697
$ws.Range("C9").Value = 0.135787664509722
$ws.Range("D9").Value = 0.0482563409154728
$ws.Range("E9").Value = 0.961558475809512
$ws.Range("F9").Value = 0.963337973541263
$ws.Range("G9").Value = ''
$ws.Range("A10").Value = 'healthy vs Grouppost_ltx:CountryNO'
$ws.Range("B10").Value = -0.358115316698667
$ws.Range("D10").Value = -2.03126987760139
$ws.Range("E10").Value = 0.0435059369349167
$ws.Range("F10").Value = 0.0978883581035626

$ws = $wb.Worksheets.Item("Simpson")
$ws.Range("A2").Value = 'healthy vs Grouppre_ltx'
$ws.Range("B2").Value = -0.0462750237387387
$ws.Range("D2").Value = -2.48087229138092
$ws.Range("E2").Value = 0.0142992845782302
$ws.Range("F2").Value = 0.128693561204072
$ws.Range("A3").Value = 'healthy , pre_ltx - CZ vs NO'
$ws.Range("B3").Value = -0.00776755151651656
$ws.Range("C3").Value = 0.0166606954098627
$ws.Range("D3").Value = -0.466220126197037
$ws.Range("E3").Value = 0.64178774498884
$ws.Range("F3").Value = 0.692428766390511
$ws.Range("A4").Value = 'healthy vs Grouppre_ltx:CountryNO'
$ws.Range("B4").Value = 0.0267112720237628
$ws.Range("D4").Value = 1.09164781225818
$ws.Range("E4").Value = 0.276876498148738
$ws.Range("F4").Value = 0.498377696667729
$ws.Range("F5").Value = 0.238294380346603
$ws.Range("F6").Value = 0.62710098994835
$ws.Range("F7").Value = 0.195148416858397
$ws.Range("A8").Value = 'healthy vs Grouppost_ltx'
$ws.Range("B8").Value = -0.00919541344462109
$ws.Range("D8").Value = -0.572018949226178
$ws.Range("E8").Value = 0.567929802357114
$ws.Range("F8").Value = 0.692428766390511
$ws.Range("A9").Value = 'healthy , post_ltx - CZ vs NO'
$ws.Range("B9").Value = -0.00776755151651663
$ws.Range("C9").Value = 0.0196093731266211
$ws.Range("D9").Value = -0.396114218764681
$ws.Range("E9").Value = 0.692428766390511
$ws.Range("F9").Value = 0.692428766390511
$ws.Range("A10").Value = 'healthy vs Grouppost_ltx:CountryNO'
$ws.Range("B10").Value = -0.0331829959998233
$ws.Range("D10").Value = -1.30333763179713
$ws.Range("E10").Value = 0.193907299646965
$ws.Range("F10").Value = 0.436291424205672

$ws = $wb.Worksheets.Item("Pielou")
$ws.Range("A2").Value = 'healthy vs Grouppre_ltx'
$ws.Range("B2").Value = -0.0404256988729569
$ws.Range("C2").Value = 0.0236644736454057
$ws.Range("D2").Value = -1.70828641611496
$ws.Range("E2").Value = 0.0898153871036228
$ws.Range("F2").Value = 0.285360193196235
$ws.Range("A3").Value = 'healthy , pre_ltx - CZ vs NO'
$ws.Range("B3").Value = -0.00833324035941611
$ws.Range("C3").Value = 0.0211372134074224
$ws.Range("D3").Value = -0.394244983896026
$ws.Range("E3").Value = 0.69400437205473
$ws.Range("F3").Value = 0.942972463000473
$ws.Range("A4").Value = 'healthy vs Grouppre_ltx:CountryNO'
$ws.Range("B4").Value = 0.0103096633488069
$ws.Range("D4").Value = 0.332106866667969
$ws.Range("E4").Value = 0.740308832328362
$ws.Range("F4").Value = 0.942972463000473
$ws.Range("F5").Value = 0.285360193196235
$ws.Range("F6").Value = 0.942972463000473
$ws.Range("F7").Value = 0.285360193196235
$ws.Range("A8").Value = 'healthy vs Grouppost_ltx'
$ws.Range("B8").Value = 0.00120233091192999
$ws.Range("D8").Value = 0.0716212699951779
$ws.Range("E8").Value = 0.942972463000473
$ws.Range("F8").Value = 0.942972463000473
$ws.Range("A9").Value = 'healthy , post_ltx - CZ vs NO'
$ws.Range("B9").Value = -0.00833324035941616
$ws.Range("C9").Value = 0.0204778737281437
$ws.Range("D9").Value = -0.406938751065908
$ws.Range("E9").Value = 0.68447343571358
$ws.Range("F9").Value = 0.942972463000473
$ws.Range("G9").Value = ''
$ws.Range("A10").Value = 'healthy vs Grouppost_ltx:CountryNO'
$ws.Range("B10").Value = -0.0405585805530275
$ws.Range("D10").Value = -1.52546742232706
$ws.Range("E10").Value = 0.128668958637096
$ws.Range("F10").Value = 0.289505156933465
